$d = $word.ActiveDocument

$d.Content.Find.Execute("650÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "315÷4=", 2) | Out-Null
$d.Content.Find.Execute("623÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "189÷7=", 2) | Out-Null
$d.Content.Find.Execute("143÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "151÷5=", 2) | Out-Null
$d.Content.Find.Execute("290÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "128÷6=", 2) | Out-Null
$d.Content.Find.Execute("674÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "211÷2=", 2) | Out-Null
$d.Content.Find.Execute("320÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "958÷8=", 2) | Out-Null
$d.Content.Find.Execute("199÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "513÷9=", 2) | Out-Null
$d.Content.Find.Execute("388÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "260÷4=", 2) | Out-Null
$d.Content.Find.Execute("301÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "520÷4=", 2) | Out-Null
$d.Content.Find.Execute("872÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "119÷3=", 2) | Out-Null
$d.Content.Find.Execute("436÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "929÷7=", 2) | Out-Null
$d.Content.Find.Execute("860÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "825÷6=", 2) | Out-Null
$d.Content.Find.Execute("702÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "698÷3=", 2) | Out-Null
$d.Content.Find.Execute("121÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "485÷2=", 2) | Out-Null
$d.Content.Find.Execute("584÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "132÷6=", 2) | Out-Null
$d.Content.Find.Execute("662÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "547÷2=", 2) | Out-Null
$d.Content.Find.Execute("604÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "100÷8=", 2) | Out-Null
$d.Content.Find.Execute("434÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "753÷8=", 2) | Out-Null
$d.Content.Find.Execute("630÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "134÷9=", 2) | Out-Null
$d.Content.Find.Execute("113÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷9=", 2) | Out-Null
$d.Content.Find.Execute("288÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "808÷3=", 2) | Out-Null
$d.Content.Find.Execute("602÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "707÷7=", 2) | Out-Null
$d.Content.Find.Execute("632÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "745÷2=", 2) | Out-Null
$d.Content.Find.Execute("397÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "489÷9=", 2) | Out-Null
$d.Content.Find.Execute("620÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=", 2) | Out-Null
